$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("医療従事者")

# Insert a new row above row 5 (the existing data rows 5.. shift down by one).
$ws.Rows.Item(5).Insert()

# New row 5: 2021-05-19 data point.
$ws.Range("A5").Value = 44335
$ws.Range("B5").Value = "(水)"
$ws.Range("C5").Formula = "=SUM(D5:E5)"
$ws.Range("D5").Value = 88163
$ws.Range("E5").Value = 179703

# Copy formatting from the row below (old row5, now row6) into new row5.
$ws.Rows.Item(6).Copy()
$ws.Rows.Item(5).PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Re-apply values/formula after format paste (paste formats only touches formatting,
# but keep values explicit to be safe).
$ws.Range("A5").Value = 44335
$ws.Range("B5").Value = "(水)"
$ws.Range("C5").Formula = "=SUM(D5:E5)"
$ws.Range("D5").Value = 88163
$ws.Range("E5").Value = 179703

# Update the totals row (row 4) for the new cumulative totals.
$ws.Range("D4").Value = 3784071
$ws.Range("E4").Value = 2214762

# Update the "as of" note text in E2.
$ws.Range("E2").Value = "（5月19日時点）"
